# Update countries & provincias Spain
# - Reorder "Groenlandia" / "Islas Malvinas" shared strings (Groenlandia now precedes Islas Malvinas)
# - Refresh "Datos actualizados" timestamp from 00:01 to 01:18
# - Update latest COVID-19 case/death counters for a set of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Groenlandia / Islas Malvinas rows (A210 <-> A211) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 01:18"

# Row 4 - Estados Unidos: B4=4245169, C4=74851, D4=2016001, E4=2080770, G4=1049, H4=148398
$ws.Range("B4").Value = 4245169
$ws.Range("C4").Value = 74851
$ws.Range("D4").Value = 2016001
$ws.Range("E4").Value = 2080770
$ws.Range("G4").Value = 1049
$ws.Range("H4").Value = 148398

# Row 5 - Brasil: B5=2348200, C5=58249, E5=670534, G5=1178, H5=85385
$ws.Range("B5").Value = 2348200
$ws.Range("C5").Value = 58249
$ws.Range("E5").Value = 670534
$ws.Range("G5").Value = 1178
$ws.Range("H5").Value = 85385

# Row 18 - Colombia: B18=233541, C18=7168, D18=113864, E18=111702, G18=287, H18=7975
$ws.Range("B18").Value = 233541
$ws.Range("C18").Value = 7168
$ws.Range("D18").Value = 113864
$ws.Range("E18").Value = 111702
$ws.Range("G18").Value = 287
$ws.Range("H18").Value = 7975

# Row 23 - Argentina: B23=153520, C23=5493, E23=85266, G23=105, H23=2807
$ws.Range("B23").Value = 153520
$ws.Range("C23").Value = 5493
$ws.Range("E23").Value = 85266
$ws.Range("G23").Value = 105
$ws.Range("H23").Value = 2807

# Row 24 - Canada: B24=113179, C24=507, D24=98837, E24=5462, G24=6, H24=8880
$ws.Range("B24").Value = 113179
$ws.Range("C24").Value = 507
$ws.Range("D24").Value = 98837
$ws.Range("E24").Value = 5462
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 8880

# Row 31 - Suecia: B31=78997, C31=82, G31=10, H31=5697
$ws.Range("B31").Value = 78997
$ws.Range("C31").Value = 82
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 5697

# Row 46 - Singapur: D46=45172, E46=4176
$ws.Range("D46").Value = 45172
$ws.Range("E46").Value = 4176

# Row 50 - Nigeria: B50=39539, C50=591, D50=16559, E50=22135, G50=12, H50=845
$ws.Range("B50").Value = 39539
$ws.Range("C50").Value = 591
$ws.Range("D50").Value = 16559
$ws.Range("E50").Value = 22135
$ws.Range("G50").Value = 12
$ws.Range("H50").Value = 845

# Row 59 - Japon: B59=27956, C59=927, D59=21328, E59=5636, G59=2, H59=992
$ws.Range("B59").Value = 27956
$ws.Range("C59").Value = 927
$ws.Range("D59").Value = 21328
$ws.Range("E59").Value = 5636
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 992

# Row 71 - Chequia: B71=15081, C71=281, D71=9422, E71=5290, G71=4, H71=369
$ws.Range("B71").Value = 15081
$ws.Range("C71").Value = 281
$ws.Range("D71").Value = 9422
$ws.Range("E71").Value = 5290
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 369

# Row 85 - Noruega: B85=9092, C85=7, E85=163
$ws.Range("B85").Value = 9092
$ws.Range("C85").Value = 7
$ws.Range("E85").Value = 163

# Row 115 - Montenegro: B115=2665, C115=96, D115=642, E115=1980, G115=3, H115=43
$ws.Range("B115").Value = 2665
$ws.Range("C115").Value = 96
$ws.Range("D115").Value = 642
$ws.Range("E115").Value = 1980
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 43

# Row 138 - Surinam: B138=1305, C138=71, D138=804, E138=478
$ws.Range("B138").Value = 1305
$ws.Range("C138").Value = 71
$ws.Range("D138").Value = 804
$ws.Range("E138").Value = 478

# Row 143 - Niger: D143=1024, E143=31
$ws.Range("D143").Value = 1024
$ws.Range("E143").Value = 31

# Row 152 - Togo: B152=839, C152=11, D152=585, E152=237, G152=1, H152=17
$ws.Range("B152").Value = 839
$ws.Range("C152").Value = 11
$ws.Range("D152").Value = 585
$ws.Range("E152").Value = 237
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 17

# Row 191 - Antigua y Barbuda: B191=82, C191=6, D191=60, E191=19
$ws.Range("B191").Value = 82
$ws.Range("C191").Value = 6
$ws.Range("D191").Value = 60
$ws.Range("E191").Value = 19

